$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "11×96="
$t.Cell(1,2).Range.Text = "84×66="
$t.Cell(1,3).Range.Text = "76×36="
$t.Cell(1,4).Range.Text = "51×41="
$t.Cell(1,5).Range.Text = "17×92="
$t.Cell(5,1).Range.Text = "75×59="
$t.Cell(5,2).Range.Text = "15×73="
$t.Cell(5,3).Range.Text = "58×45="
$t.Cell(5,4).Range.Text = "49×53="
$t.Cell(5,5).Range.Text = "24×97="
$t.Cell(10,1).Range.Text = "44×81="
$t.Cell(10,2).Range.Text = "46×48="
$t.Cell(10,3).Range.Text = "14×19="
$t.Cell(10,4).Range.Text = "89×38="
$t.Cell(10,5).Range.Text = "15×68="
$t.Cell(15,1).Range.Text = "39×90="
$t.Cell(15,2).Range.Text = "40×60="
$t.Cell(15,3).Range.Text = "43×32="
$t.Cell(15,4).Range.Text = "25×67="
$t.Cell(15,5).Range.Text = "99×97="
$t.Cell(20,1).Range.Text = "20×79="
$t.Cell(20,2).Range.Text = "49×93="
$t.Cell(20,3).Range.Text = "53×17="
$t.Cell(20,4).Range.Text = "72×18="
$t.Cell(20,5).Range.Text = "34×67="
